$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 894.06665
$ws.Range("J17").Value = 894.06665
$ws.Range("L17").Value = 2682.19995
$ws.Range("N17").Value = -3018.19995
$ws.Range("H19").Value = 21547258
$ws.Range("I19").Value = 15653054
$ws.Range("J19").Value = 33335668
$ws.Range("K19").Value = 15653054
$ws.Range("L19").Value = 33335668
$ws.Range("M19").Value = -15652879
$ws.Range("N19").Value = -33336018
$ws.Range("H32").Value = 1588.4166
$ws.Range("I32").Value = 757.375
$ws.Range("J32").Value = 3250.5
$ws.Range("K32").Value = 757.375
$ws.Range("L32").Value = 3250.5
$ws.Range("M32").Value = -431.375
$ws.Range("N32").Value = -3902.5
$ws.Range("H98").Value = 3750
$ws.Range("I98").Value = 2500
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 2500
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -1002
$ws.Range("N98").Value = -7996
$ws.Range("H117").Value = 41801
$ws.Range("J117").Value = 41801
$ws.Range("L117").Value = 41801
$ws.Range("N117").Value = -50979
$ws.Range("H122").Value = 3750
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -19900
$ws.Range("H137").Value = 12552515
$ws.Range("I137").Value = 50002920
$ws.Range("J137").Value = 69046.92999999999
$ws.Range("K137").Value = 150008760
$ws.Range("L137").Value = 207140.79
$ws.Range("M137").Value = -150006210
$ws.Range("N137").Value = -212240.79
$ws.Range("H140").Value = 54663.332
$ws.Range("J140").Value = 54663.332
$ws.Range("L140").Value = 54663.332
$ws.Range("N140").Value = -65023.332

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2393.0588
$ws.Range("I2").Value = 2123
$ws.Range("K2").Value = 2123
$ws.Range("M2").Value = -2010
$ws.Range("H61").Value = 3306.8462
$ws.Range("I61").Value = 1689.3334
$ws.Range("J61").Value = 4693.2856
$ws.Range("K61").Value = 1689.3334
$ws.Range("L61").Value = 4693.2856
$ws.Range("M61").Value = -1477.3334
$ws.Range("N61").Value = -5117.2856
$ws.Range("H74").Value = 7337.1177
$ws.Range("I74").Value = 7854.467
$ws.Range("J74").Value = 3457
$ws.Range("K74").Value = 7854.467
$ws.Range("L74").Value = 3457
$ws.Range("M74").Value = -6980.467
$ws.Range("N74").Value = -5205
$ws.Range("H77").Value = 7337.1177
$ws.Range("I77").Value = 7854.467
$ws.Range("J77").Value = 3457
$ws.Range("K77").Value = 39272.335
$ws.Range("L77").Value = 17285
$ws.Range("M77").Value = -34904.335
$ws.Range("N77").Value = -26021
$ws.Range("H92").Value = 24699.666
$ws.Range("J92").Value = 24699.666
$ws.Range("L92").Value = 24699.666
$ws.Range("N92").Value = -29691.666
$ws.Range("H116").Value = 2393.0588
$ws.Range("I116").Value = 2123
$ws.Range("K116").Value = 2123
$ws.Range("M116").Value = 171
$ws.Range("H132").Value = 2041.3462
$ws.Range("I132").Value = 1119.1578
$ws.Range("J132").Value = 4544.4287
$ws.Range("K132").Value = 3357.4734
$ws.Range("L132").Value = 13633.2861
$ws.Range("M132").Value = -827.4733999999999
$ws.Range("N132").Value = -18693.2861
$ws.Range("H136").Value = 3306.8462
$ws.Range("I136").Value = 1689.3334
$ws.Range("J136").Value = 4693.2856
$ws.Range("K136").Value = 5068.0002
$ws.Range("L136").Value = 14079.8568
$ws.Range("M136").Value = -2518.0002
$ws.Range("N136").Value = -19179.8568

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2393.0588
$ws.Range("I3").Value = 2123
$ws.Range("K3").Value = 2123
$ws.Range("M3").Value = -2009
$ws.Range("H20").Value = 3871.353
$ws.Range("I20").Value = 2711.3
$ws.Range("K20").Value = 2711.3
$ws.Range("M20").Value = -2464.3
$ws.Range("H99").Value = 1759.7142
$ws.Range("I99").Value = 465.75
$ws.Range("J99").Value = 3485
$ws.Range("K99").Value = 465.75
$ws.Range("L99").Value = 3485
$ws.Range("M99").Value = 1032.25
$ws.Range("N99").Value = -6481
$ws.Range("H134").Value = 1202.1212
$ws.Range("I134").Value = 892.129
$ws.Range("J134").Value = 6007
$ws.Range("K134").Value = 2676.387
$ws.Range("L134").Value = 18021
$ws.Range("M134").Value = -141.3870000000002
$ws.Range("N134").Value = -23091

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2694.0356
$ws.Range("I31").Value = 1733.5714
$ws.Range("J31").Value = 5575.4287
$ws.Range("K31").Value = 1733.5714
$ws.Range("L31").Value = 5575.4287
$ws.Range("M31").Value = -1438.5714
$ws.Range("N31").Value = -6165.4287
$ws.Range("H34").Value = 2694.0356
$ws.Range("I34").Value = 1733.5714
$ws.Range("J34").Value = 5575.4287
$ws.Range("K34").Value = 1733.5714
$ws.Range("L34").Value = 5575.4287
$ws.Range("M34").Value = -1531.5714
$ws.Range("N34").Value = -5979.4287
$ws.Range("H58").Value = 2810.8667
$ws.Range("I58").Value = 2364.9
$ws.Range("J58").Value = 3702.8
$ws.Range("K58").Value = 2364.9
$ws.Range("L58").Value = 3702.8
$ws.Range("M58").Value = -2161.9
$ws.Range("N58").Value = -4108.8
$ws.Range("H99").Value = 5000
$ws.Range("J99").Value = 5000
$ws.Range("L99").Value = 5000
$ws.Range("N99").Value = -7996
$ws.Range("H126").Value = 5000
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2586.3333
$ws.Range("I132").Value = 1445.1818
$ws.Range("K132").Value = 4335.5454
$ws.Range("M132").Value = -1805.5454
$ws.Range("H134").Value = 2037
$ws.Range("I134").Value = 1287.3043
$ws.Range("K134").Value = 3861.9129
$ws.Range("M134").Value = -1326.9129
$ws.Range("H136").Value = 2810.8667
$ws.Range("I136").Value = 2364.9
$ws.Range("J136").Value = 3702.8
$ws.Range("K136").Value = 7094.700000000001
$ws.Range("L136").Value = 11108.4
$ws.Range("M136").Value = -4544.700000000001
$ws.Range("N136").Value = -16208.4

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 15.384615
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 16.5
$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 49.5
$ws.Range("M12").Value = 167
$ws.Range("N12").Value = -395.5
$ws.Range("H109").Value = 3606.1482
$ws.Range("I109").Value = 1665.875
$ws.Range("K109").Value = 4997.625
$ws.Range("M109").Value = -3957.625
$ws.Range("H131").Value = 785.45
$ws.Range("I131").Value = 562.25
$ws.Range("J131").Value = 934.25
$ws.Range("K131").Value = 1686.75
$ws.Range("L131").Value = 2802.75
$ws.Range("M131").Value = 3353.25
$ws.Range("N131").Value = -12882.75

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2703.6
$ws.Range("I102").Value = 2570.3333
$ws.Range("K102").Value = 2570.3333
$ws.Range("M102").Value = -948.3332999999998
$ws.Range("H116").Value = 57907.332
$ws.Range("J116").Value = 57907.332
$ws.Range("L116").Value = 57907.332
$ws.Range("N116").Value = -67085.33199999999
$ws.Range("H132").Value = 6610.7144
$ws.Range("I132").Value = 6951
$ws.Range("K132").Value = 20853
$ws.Range("M132").Value = -18323

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2514.6
$ws.Range("I132").Value = 2249.577
$ws.Range("J132").Value = 4237.25
$ws.Range("K132").Value = 6748.731000000001
$ws.Range("L132").Value = 12711.75
$ws.Range("M132").Value = -4218.731000000001
$ws.Range("N132").Value = -17771.75

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2334.5386
$ws.Range("I132").Value = 1901.2593
$ws.Range("J132").Value = 3309.4167
$ws.Range("K132").Value = 5703.7779
$ws.Range("L132").Value = 9928.250100000001
$ws.Range("M132").Value = -3173.7779
$ws.Range("N132").Value = -14988.2501
$ws.Range("H136").Value = 4766.4287
$ws.Range("I136").Value = 5045.4346
$ws.Range("K136").Value = 15136.3038
$ws.Range("M136").Value = -12586.3038

Write-Host "Applied all changes"